$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-23 15:23:56"

$wsZhCn.Range("H2").Value = "2016-08-23 15:23:50"
$wsZhCn.Range("K2").Value = "2016-08-23 15:24:27"

$wsDeDe.Range("H2").Value = "2016-08-23 15:23:56"
$wsDeDe.Range("K2").Value = "2016-08-23 15:24:35"
